$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.005.43'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '2.297.81'
$ws.Range("E4").Value = '  -0.01%  '
$cell = $ws.Range("D5")
$cell.Value = "'300.73"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'99.41"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.81%  '
$cell = $ws.Range("D7")
$cell.Value = "'0.509"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  -0.04%  '
$cell = $ws.Range("D9")
$cell.Value = "'0.509"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.34%  '
$cell = $ws.Range("D10")
$cell.Value = "'36.13"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +6.90%  '
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("E12").Value = '  +0.58%  '
$cell = $ws.Range("D13")
$cell.Value = "'17.58"
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("E14").Value = '  +1.50%  '
$ws.Range("D15").Value = '2.654.44'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '2.292.46'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '42.904.26'
$ws.Range("E18").Value = '  -0.31%  '
$cell = $ws.Range("D19")
$cell.Value = "'12.80"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +8.86%  '
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("E22").Value = '  +0.24%  '
$cell = $ws.Range("D23")
$cell.Value = "'235.33"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '
$cell = $ws.Range("D24")
$cell.Value = "'2.18"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +6.60%  '
$cell = $ws.Range("D25")
$cell.Value = "'1.01"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.53%  '
$cell = $ws.Range("D26")
$cell.Value = "'2.45"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.64%  '
$cell = $ws.Range("D27")
$cell.Value = "'24.91"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.83%  '
$cell = $ws.Range("D28")
$cell.Value = "'169.55"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.71%  '
$cell = $ws.Range("D29")
$cell.Value = "'34.39"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("E30").Value = '  -5.46%  '
$ws.Range("E31").Value = '  +0.36%  '
$cell = $ws.Range("D32")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  +1.71%  '
$cell = $ws.Range("D34")
$cell.Value = "'17.62"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.82%  '
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("D43").Value = '1.986.89'
$ws.Range("E43").Value = '  +0.19%  '
$cell = $ws.Range("D44")
$cell.Value = "'2.25"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -5.66%  '
$cell = $ws.Range("D45")
$cell.Value = "'10.14"
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("E46").Value = '  +1.05%  '
$cell = $ws.Range("D47")
$cell.Value = "'17.49"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.39%  '
$cell = $ws.Range("D48")
$cell.Value = "'55.41"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +3.86%  '
$ws.Range("E49").Value = '  +3.58%  '
$ws.Range("D50").Value = '2.522.06'
$ws.Range("E50").Value = '  -0.69%  '
$cell = $ws.Range("D51")
$cell.Value = "'70.68"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
